$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(4)
$tbl = $sh.Table
$cell = $tbl.Cell(1, 1)
$cell.Shape.TextFrame.TextRange.Text = ".NET 6"
